$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 10: action value changes from "oauthclick" to "clickloc"
$ws.Range("A10").Value = "clickloc"

# Remove rows 11 and 12 entirely (cells + formatting), while keeping row numbers of
# subsequent rows (15, 20) unchanged - i.e. clear, not delete-and-shift.
$ws.Range("A11:E12").Clear()

# Reset the selection to the default top-left cell.
$ws.Range("A1").Select()
